$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 109.7  # was 414.45456
$ws.Range("I41").Value = 110.77778  # was 114.57143
$ws.Range("J41").Value = 100  # was 939.25
$ws.Range("K41").Value = 110.77778  # was 114.57143
$ws.Range("L41").Value = 100  # was 939.25
$ws.Range("M41").Value = 329.22222  # was 325.42857
$ws.Range("N41").Value = -980  # was -1819.25
$ws.Range("H64").Value = 3767.0527  # was 3799.389
$ws.Range("I64").Value = 3459.5386  # was 3482.4167
$ws.Range("K64").Value = 3459.5386  # was 3482.4167
$ws.Range("M64").Value = -3211.5386  # was -3234.4167
$ws.Range("H67").Value = 3767.0527  # was 3799.389
$ws.Range("I67").Value = 3459.5386  # was 3482.4167
$ws.Range("K67").Value = 3459.5386  # was 3482.4167
$ws.Range("M67").Value = -2601.5386  # was -2624.4167
$ws.Range("H74").Value = 5637.5835  # was 5338
$ws.Range("I74").Value = 4381.857  # was 4672.5
$ws.Range("J74").Value = 7395.6  # was 8000
$ws.Range("K74").Value = 4381.857  # was 4672.5
$ws.Range("L74").Value = 7395.6  # was 8000
$ws.Range("M74").Value = -3445.857  # was -3736.5
$ws.Range("N74").Value = -9267.6  # was -9872
$ws.Range("H76").Value = 4016.3022  # was 4043.0977
$ws.Range("I76").Value = 3834.3142  # was 3840.75
$ws.Range("J76").Value = 4812.5  # was 5500
$ws.Range("K76").Value = 3834.3142  # was 3840.75
$ws.Range("L76").Value = 4812.5  # was 5500
$ws.Range("M76").Value = -3519.3142  # was -3525.75
$ws.Range("N76").Value = -5442.5  # was -6130
$ws.Range("H77").Value = 5637.5835  # was 5338
$ws.Range("I77").Value = 4381.857  # was 4672.5
$ws.Range("J77").Value = 7395.6  # was 8000
$ws.Range("K77").Value = 21909.285  # was 23362.5
$ws.Range("L77").Value = 36978  # was 40000
$ws.Range("M77").Value = -17229.285  # was -18682.5
$ws.Range("N77").Value = -46338  # was -49360
$ws.Range("H79").Value = 4016.3022  # was 4043.0977
$ws.Range("I79").Value = 3834.3142  # was 3840.75
$ws.Range("J79").Value = 4812.5  # was 5500
$ws.Range("K79").Value = 3834.3142  # was 3840.75
$ws.Range("L79").Value = 4812.5  # was 5500
$ws.Range("M79").Value = -2742.3142  # was -2748.75
$ws.Range("N79").Value = -6996.5  # was -7684
$ws.Range("H98").Value = 2850.8948  # was 3072.5625
$ws.Range("I98").Value = 2162.9167  # was 2395.5
$ws.Range("J98").Value = 4030.2856  # was 4201
$ws.Range("K98").Value = 2162.9167  # was 2395.5
$ws.Range("L98").Value = 4030.2856  # was 4201
$ws.Range("M98").Value = -664.9167000000002  # was -897.5
$ws.Range("N98").Value = -7026.2856  # was -7197
$ws.Range("H113").Value = 2750.1765  # was 2856.4666
$ws.Range("J113").Value = 3011  # was 3203.3635
$ws.Range("L113").Value = 3011  # was 3203.3635
$ws.Range("N113").Value = -9519  # was -9711.363499999999
$ws.Range("H122").Value = 2850.8948  # was 3072.5625
$ws.Range("I122").Value = 2162.9167  # was 2395.5
$ws.Range("J122").Value = 4030.2856  # was 4201
$ws.Range("K122").Value = 6488.750100000001  # was 7186.5
$ws.Range("L122").Value = 12090.8568  # was 12603
$ws.Range("M122").Value = -4038.750100000001  # was -4736.5
$ws.Range("N122").Value = -16990.8568  # was -17503
$ws.Range("H129").Value = 935.9861  # was 936.775
$ws.Range("I129").Value = 422.1111  # was 505.7
$ws.Range("J129").Value = 1009.39685  # was 998.3570999999999
$ws.Range("K129").Value = 1266.3333  # was 1517.1
$ws.Range("L129").Value = 3028.19055  # was 2995.0713
$ws.Range("M129").Value = 3733.6667  # was 3482.9
$ws.Range("N129").Value = -13028.19055  # was -12995.0713
$ws.Range("H132").Value = 1499.4849  # was 1657.7241
$ws.Range("I132").Value = 1499.4849  # was 1694.2142
$ws.Range("J132").Value = 0  # was 636
$ws.Range("K132").Value = 4498.4547  # was 5082.642599999999
$ws.Range("L132").Value = 0  # was 1908
$ws.Range("M132").Value = -1968.4547  # was -2552.642599999999
$ws.Range("N132").ClearContents()  # was -6968

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 38249  # was 40249
$ws.Range("J62").Value = 38249  # was 40249
$ws.Range("L62").Value = 38249  # was 40249
$ws.Range("N62").Value = -39497  # was -41497
$ws.Range("H63").Value = 4191.25  # was 2850.625
$ws.Range("I63").Value = 3265  # was 2179.8
$ws.Range("J63").Value = 4500  # was 3968.6667
$ws.Range("K63").Value = 3265  # was 2179.8
$ws.Range("L63").Value = 4500  # was 3968.6667
$ws.Range("M63").Value = -2579  # was -1493.8
$ws.Range("N63").Value = -5872  # was -5340.6667
$ws.Range("H65").Value = 38249  # was 40249
$ws.Range("J65").Value = 38249  # was 40249
$ws.Range("L65").Value = 114747  # was 120747
$ws.Range("N65").Value = -120987  # was -126987
$ws.Range("H66").Value = 4191.25  # was 2850.625
$ws.Range("I66").Value = 3265  # was 2179.8
$ws.Range("J66").Value = 4500  # was 3968.6667
$ws.Range("K66").Value = 16325  # was 10899
$ws.Range("L66").Value = 22500  # was 19843.3335
$ws.Range("M66").Value = -12893  # was -7467
$ws.Range("N66").Value = -29364  # was -26707.3335
$ws.Range("H74").Value = 1692.475  # was 1799.919
$ws.Range("I74").Value = 1678.4722  # was 1759.5
$ws.Range("J74").Value = 1818.5  # was 2258
$ws.Range("K74").Value = 1678.4722  # was 1759.5
$ws.Range("L74").Value = 1818.5  # was 2258
$ws.Range("M74").Value = -804.4721999999999  # was -885.5
$ws.Range("N74").Value = -3566.5  # was -4006
$ws.Range("H77").Value = 1692.475  # was 1799.919
$ws.Range("I77").Value = 1678.4722  # was 1759.5
$ws.Range("J77").Value = 1818.5  # was 2258
$ws.Range("K77").Value = 8392.360999999999  # was 8797.5
$ws.Range("L77").Value = 9092.5  # was 11290
$ws.Range("M77").Value = -4024.360999999999  # was -4429.5
$ws.Range("N77").Value = -17828.5  # was -20026
$ws.Range("H122").Value = 1639.9231  # was 1758.091
$ws.Range("I122").Value = 1684.9166  # was 1758.091
$ws.Range("J122").Value = 1100  # was 0
$ws.Range("K122").Value = 5054.7498  # was 5274.272999999999
$ws.Range("L122").Value = 3300  # was 0
$ws.Range("M122").Value = -2604.7498  # was -2824.272999999999
$ws.Range("N122").Value = -8200  # new cell
$ws.Range("H132").Value = 2285.6  # was 2214.0435
$ws.Range("I132").Value = 1955.8182  # was 1795
$ws.Range("J132").Value = 2688.6667  # was 2999.75
$ws.Range("K132").Value = 5867.4546  # was 5385
$ws.Range("L132").Value = 8066.000100000001  # was 8999.25
$ws.Range("M132").Value = -3337.4546  # was -2855
$ws.Range("N132").Value = -13126.0001  # was -14059.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2271.875  # was 2237.675
$ws.Range("I31").Value = 1433.35  # was 1380.8572
$ws.Range("J31").Value = 3110.4  # was 3184.6843
$ws.Range("K31").Value = 1433.35  # was 1380.8572
$ws.Range("L31").Value = 3110.4  # was 3184.6843
$ws.Range("M31").Value = -1138.35  # was -1085.8572
$ws.Range("N31").Value = -3700.4  # was -3774.6843
$ws.Range("H34").Value = 2271.875  # was 2237.675
$ws.Range("I34").Value = 1433.35  # was 1380.8572
$ws.Range("J34").Value = 3110.4  # was 3184.6843
$ws.Range("K34").Value = 1433.35  # was 1380.8572
$ws.Range("L34").Value = 3110.4  # was 3184.6843
$ws.Range("M34").Value = -1231.35  # was -1178.8572
$ws.Range("N34").Value = -3514.4  # was -3588.6843
$ws.Range("H62").Value = 2875  # was 3000
$ws.Range("I62").Value = 2875  # was 3000
$ws.Range("K62").Value = 2875  # was 3000
$ws.Range("M62").Value = -2251  # was -2376
$ws.Range("H65").Value = 2875  # was 3000
$ws.Range("I65").Value = 2875  # was 3000
$ws.Range("K65").Value = 14375  # was 15000
$ws.Range("M65").Value = -11255  # was -11880
$ws.Range("H99").Value = 1200.2307  # was 1170.8572
$ws.Range("I99").Value = 1091.1818  # was 1066
$ws.Range("K99").Value = 1091.1818  # was 1066
$ws.Range("M99").Value = 406.8181999999999  # was 432
$ws.Range("H126").Value = 1200.2307  # was 1170.8572
$ws.Range("I126").Value = 1091.1818  # was 1066
$ws.Range("K126").Value = 3273.5454  # was 3198
$ws.Range("M126").Value = -803.5454  # was -728
$ws.Range("H132").Value = 2560.3416  # was 2328.8298
$ws.Range("I132").Value = 2186.75  # was 1959.3948
$ws.Range("K132").Value = 6560.25  # was 5878.1844
$ws.Range("M132").Value = -4030.25  # was -3348.1844
$ws.Range("H134").Value = 2483.7837  # was 2308.1462
$ws.Range("I134").Value = 2241.7666  # was 2096.1516
$ws.Range("J134").Value = 3521  # was 3182.625
$ws.Range("K134").Value = 6725.2998  # was 6288.4548
$ws.Range("L134").Value = 10563  # was 9547.875
$ws.Range("M134").Value = -4190.2998  # was -3753.4548
$ws.Range("N134").Value = -15633  # was -14617.875

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 18530186  # was 23824296
$ws.Range("I5").Value = 599  # was 464
$ws.Range("J5").Value = 33353856  # was 41692170
$ws.Range("K5").Value = 1797  # was 1392
$ws.Range("L5").Value = 100061568  # was 125076510
$ws.Range("M5").Value = -1685  # was -1280
$ws.Range("N5").Value = -100061792  # was -125076734
$ws.Range("H132").Value = 1346.5714  # was 1405.8334
$ws.Range("I132").Value = 994.7143  # was 996.6667
$ws.Range("J132").Value = 1698.4286  # was 1815
$ws.Range("K132").Value = 8952.4287  # was 8970.0003
$ws.Range("L132").Value = 15285.8574  # was 16335
$ws.Range("M132").Value = -6422.4287  # was -6440.0003
$ws.Range("N132").Value = -20345.8574  # was -21395
$ws.Range("H133").Value = 4468.2856  # was 3699.1667
$ws.Range("I133").Value = 4630  # was 2870.1667
$ws.Range("J133").Value = 4347  # was 4113.6665
$ws.Range("K133").Value = 13890  # was 8610.500100000001
$ws.Range("L133").Value = 13041  # was 12340.9995
$ws.Range("M133").Value = -8830  # was -3550.500100000001
$ws.Range("N133").Value = -23161  # was -22460.9995
$ws.Range("H135").Value = 18530186  # was 23824296
$ws.Range("I135").Value = 599  # was 464
$ws.Range("J135").Value = 33353856  # was 41692170
$ws.Range("K135").Value = 5391  # was 4176
$ws.Range("L135").Value = 300184704  # was 375229530
$ws.Range("M135").Value = -2856  # was -1641
$ws.Range("N135").Value = -300189774  # was -375234600

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5736.353  # was 5764.592
$ws.Range("I70").Value = 5553.1514  # was 5588.533
$ws.Range("J70").Value = 6072.222  # was 6042.579
$ws.Range("K70").Value = 5553.1514  # was 5588.533
$ws.Range("L70").Value = 6072.222  # was 6042.579
$ws.Range("M70").Value = -5283.1514  # was -5318.533
$ws.Range("N70").Value = -6612.222  # was -6582.579
$ws.Range("H73").Value = 5736.353  # was 5764.592
$ws.Range("I73").Value = 5553.1514  # was 5588.533
$ws.Range("J73").Value = 6072.222  # was 6042.579
$ws.Range("K73").Value = 5553.1514  # was 5588.533
$ws.Range("L73").Value = 6072.222  # was 6042.579
$ws.Range("M73").Value = -4617.1514  # was -4652.533
$ws.Range("N73").Value = -7944.222  # was -7914.579
$ws.Range("H80").Value = 7150.278  # was 8013.3335
$ws.Range("I80").Value = 12161  # was 12240
$ws.Range("J80").Value = 5223.077  # was 5900
$ws.Range("K80").Value = 12161  # was 12240
$ws.Range("L80").Value = 5223.077  # was 5900
$ws.Range("M80").Value = -11163  # was -11242
$ws.Range("N80").Value = -7219.077  # was -7896
$ws.Range("H83").Value = 7150.278  # was 8013.3335
$ws.Range("I83").Value = 12161  # was 12240
$ws.Range("J83").Value = 5223.077  # was 5900
$ws.Range("K83").Value = 60805  # was 61200
$ws.Range("L83").Value = 26115.385  # was 29500
$ws.Range("M83").Value = -55813  # was -56208
$ws.Range("N83").Value = -36099.385  # was -39484
$ws.Range("H116").Value = 79800  # was 78500
$ws.Range("J116").Value = 79800  # was 78500
$ws.Range("L116").Value = 79800  # was 78500
$ws.Range("N116").Value = -88978  # was -87678
$ws.Range("H122").Value = 6209.8  # was 8556.857
$ws.Range("I122").Value = 7451.625  # was 9735.5
$ws.Range("J122").Value = 1242.5  # was 1485
$ws.Range("K122").Value = 22354.875  # was 29206.5
$ws.Range("L122").Value = 3727.5  # was 4455
$ws.Range("M122").Value = -19904.875  # was -26756.5
$ws.Range("N122").Value = -8627.5  # was -9355
$ws.Range("H132").Value = 3552.75  # was 8798.333000000001
$ws.Range("I132").Value = 3587.0833  # was 3427.4614
$ws.Range("J132").Value = 3449.75  # was 22762.6
$ws.Range("K132").Value = 10761.2499  # was 10282.3842
$ws.Range("L132").Value = 10349.25  # was 68287.79999999999
$ws.Range("M132").Value = -8231.249899999999  # was -7752.3842
$ws.Range("N132").Value = -15409.25  # was -73347.79999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3745.9834  # was 4156.393
$ws.Range("I136").Value = 2259.8438  # was 2663.3
$ws.Range("J136").Value = 5444.4287  # was 5879.1924
$ws.Range("K136").Value = 6779.5314  # was 7989.900000000001
$ws.Range("L136").Value = 16333.2861  # was 17637.5772
$ws.Range("M136").Value = -4229.5314  # was -5439.900000000001
$ws.Range("N136").Value = -21433.2861  # was -22737.5772

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2485  # was 2787.9512
$ws.Range("I136").Value = 1567.742  # was 1811.5769
$ws.Range("J136").Value = 4262.1875  # was 4480.3335
$ws.Range("K136").Value = 4703.226  # was 5434.7307
$ws.Range("L136").Value = 12786.5625  # was 13441.0005
$ws.Range("M136").Value = -2153.226  # was -2884.7307
$ws.Range("N136").Value = -17886.5625  # was -18541.0005
